# Refresh the cryptocurrency price/volume table with newly scraped values.
#
# Column D ("Price") holds text that sometimes *looks* numeric (e.g.
# "316.35"); Excel would silently reinterpret such a literal as a number
# when assigned to a General-formatted cell, so for those values we force
# the cell to Text format ("@") before assigning, which is what keeps the
# value stored as a string (matching the source data's formatting, e.g.
# trailing zeros like "0.3894" or "122.01" that a numeric type would not
# preserve). Values that already contain multiple dots (e.g. "28.513.17")
# aren't valid numbers anyway, so Excel keeps them as text automatically.
#
# Column E ("Volume(1h)") values contain surrounding spaces and a percent
# sign, so they are always stored as text without any extra handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "28.513.17";  E = "  +1.23%  " }
    @{ Row = 3;  D = "1.876.00";   E = "  +0.83%  " }
    @{ Row = 4;  D = $null;        E = "  +0.27%  " }
    @{ Row = 5;  D = "316.35";     E = "  +1.18%  " }
    @{ Row = 6;  D = $null;        E = "  +0.64%  " }
    @{ Row = 7;  D = "0.5082";     E = "  -0.21%  " }
    @{ Row = 8;  D = "0.3894";     E = "  -0.16%  " }
    @{ Row = 9;  D = "0.08389";    E = "  +1.52%  " }
    @{ Row = 10; D = "1.103";      E = "  -0.63%  " }
    @{ Row = 11; D = "41.79";      E = "  +0.45%  " }
    @{ Row = 12; D = "6.221";      E = "  +0.09%  " }
    @{ Row = 13; D = "1.873.48";   E = "  +0.98%  " }
    @{ Row = 14; D = "20.38";      E = "  +0.83%  " }
    @{ Row = 15; D = "7.237";      E = "  +0.62%  " }
    @{ Row = 16; D = "1.011";      E = "  +0.43%  " }
    @{ Row = 17; D = "0.00001103"; E = "  +0.62%  " }
    @{ Row = 18; D = "91.36";      E = "  +0.49%  " }
    @{ Row = 19; D = "0.06731";    E = "  +1.00%  " }
    @{ Row = 20; D = "17.70";      E = "  +0.62%  " }
    @{ Row = 21; D = "1.009";      E = "  +0.48%  " }
    @{ Row = 22; D = "5.931";      E = "  -0.01%  " }
    @{ Row = 23; D = "28.558.90";  E = "  +1.38%  " }
    @{ Row = 24; D = "11.08";      E = "  +0.28%  " }
    @{ Row = 25; D = "2.234";      E = "  +0.30%  " }
    @{ Row = 26; D = "2.088.29";   E = "  +2.02%  " }
    @{ Row = 27; D = "161.66";     E = "  +1.32%  " }
    @{ Row = 28; D = "20.63";      E = "  +0.77%  " }
    @{ Row = 29; D = $null;        E = "  -1.00%  " }
    @{ Row = 30; D = "125.89";     E = "  +0.56%  " }
    @{ Row = 31; D = "0.1045";     E = "  -0.72%  " }
    @{ Row = 32; D = "1.038";      E = "  +0.53%  " }
    @{ Row = 33; D = "5.777";      E = "  -0.66%  " }
    @{ Row = 34; D = "3.616";      E = "  +0.44%  " }
    @{ Row = 35; D = "0.02463";    E = "  +1.40%  " }
    @{ Row = 36; D = "0.06541";    E = "  +1.30%  " }
    @{ Row = 37; D = "0.2161";     E = "  -0.19%  " }
    @{ Row = 38; D = "8.861";      E = "  -2.02%  " }
    @{ Row = 39; D = "5.072";      E = "  +2.71%  " }
    @{ Row = 40; D = "1.251";      E = "  +0.43%  " }
    @{ Row = 41; D = "1.190";      E = "  +0.97%  " }
    @{ Row = 42; D = "0.6405";     E = "  -0.09%  " }
    @{ Row = 43; D = "11.10";      E = "  +0.59%  " }
    @{ Row = 44; D = $null;        E = "  +0.70%  " }
    @{ Row = 45; D = "0.6023";     E = "  +0.39%  " }
    @{ Row = 46; D = "12.98";      E = "  +0.06%  " }
    @{ Row = 47; D = "3.695";      E = "  +0.94%  " }
    @{ Row = 48; D = "2.006";      E = "  +0.52%  " }
    @{ Row = 49; D = "1.215";      E = "  +1.09%  " }
    @{ Row = 50; D = "122.01";     E = "  +1.10%  " }
    @{ Row = 51; D = $null;        E = "  -10.87%  " }
)

foreach ($u in $updates) {
    $r = $u.Row

    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($r, 4)
        # Only values that would parse as a plain number need the Text
        # format nudge; the "xx.xxx.xx" thousand-grouped prices already
        # fail numeric parsing and stay text on their own.
        $looksNumeric = $u.D -match '^\d+(\.\d+)?$'
        if ($looksNumeric) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }

    $ws.Cells.Item($r, 5).Value = $u.E
}
